$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "git@github.com:Nowpngs/universal-feature-update.git"
$ws.Range("C2").Value = "9a2ae80c1ffa646c0829324b8bc02dfa5aa799c8"

$ws.Range("A2").Select()
